$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.589.21'
$ws.Range("E2").Value = '  +3.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.606.30'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.38'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("E6").Value = '  +2.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.76'
$ws.Range("E8").Value = '  +7.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.54'
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("E10").Value = '  +2.31%  '
$ws.Range("E11").Value = '  +2.40%  '
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.835.81'
$ws.Range("E13").Value = '  +2.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.628.01'
$ws.Range("E14").Value = '  +4.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.592.77'
$ws.Range("E15").Value = '  +3.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.537'
$ws.Range("E16").Value = '  +3.98%  '
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.43'
$ws.Range("E18").Value = '  +3.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.73'
$ws.Range("E19").Value = '  +5.68%  '
$ws.Range("E20").Value = '  +3.97%  '
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.22'
$ws.Range("E24").Value = '  +2.23%  '
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.59'
$ws.Range("E26").Value = '  +1.95%  '
$ws.Range("E27").Value = '  +2.82%  '
$ws.Range("E28").Value = '  +3.47%  '
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +3.28%  '
$ws.Range("E33").Value = '  +1.51%  '
$ws.Range("E34").Value = '  +4.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.406.70'
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("E37").Value = '  +5.12%  '
$ws.Range("E38").Value = '  +5.40%  '
$ws.Range("E39").Value = '  +0.20%  '
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("E41").Value = '  +4.16%  '
$ws.Range("E42").Value = '  +2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0492'
$ws.Range("E43").Value = '  +6.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.07'
$ws.Range("E44").Value = '  +27.51%  '
$ws.Range("E45").Value = '  +4.02%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.97'
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.746.02'
$ws.Range("E49").Value = '  +2.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.864'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.65'
$ws.Range("E51").Value = '  +2.23%  '
